# edit.ps1 - "edit protocol and data sheets"
#
# Fills in the previously-blank syllable-count columns (J/K/L/M/N,
# S/T/U, Z/AA/AB) for rows 11-13 of Sheet1, then moves the active
# selection down to A15 (past the data that was just entered).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 11 (SRE26M04, song 6) ---------------------------------------
$ws.Range("J11").Value  = 3
$ws.Range("K11").Value  = 5
$ws.Range("L11").Value  = 1
$ws.Range("M11").Value  = 1
$ws.Range("N11").Value  = 2
$ws.Range("S11").Value  = 1
$ws.Range("T11").Value  = 1
$ws.Range("U11").Value  = 1
$ws.Range("Z11").Value  = 1
$ws.Range("AA11").Value = 1
$ws.Range("AB11").Value = 2

# --- Row 12 (SRE26M04, song 4) ---------------------------------------
$ws.Range("J12").Value  = 2
$ws.Range("K12").Value  = 5
$ws.Range("L12").Value  = 1
$ws.Range("M12").Value  = 1
$ws.Range("N12").Value  = 2
$ws.Range("S12").Value  = 1
$ws.Range("T12").Value  = 1
$ws.Range("U12").Value  = 3

# --- Row 13 (SRE26M04, song 2) ---------------------------------------
$ws.Range("J13").Value  = 3
$ws.Range("K13").Value  = 5
$ws.Range("L13").Value  = 1
$ws.Range("M13").Value  = 1
$ws.Range("N13").Value  = 2
$ws.Range("S13").Value  = 1
$ws.Range("T13").Value  = 1
$ws.Range("U13").Value  = 1
$ws.Range("Z13").Value  = 1
$ws.Range("AA13").Value = 1
$ws.Range("AB13").Value = 2

# Move the selection off the data that was just filled in, matching
# where the cursor ends up after the last entry (A15).
$ws.Range("A15").Select()
